$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns G, I, K hold values that look like numbers/dates ("44445",
# "01/11/25", "15/11/25") but must stay as plain text, matching the rest
# of the sheet (every other data cell is also stored as text). Force a
# text number-format before assigning so Excel doesn't auto-convert them,
# then drop the temporary formatting so the new row's cells end up with
# the sheet's default (unstyled) look, same as existing rows.
$ws.Range("G11").NumberFormat = "@"
$ws.Range("I11").NumberFormat = "@"
$ws.Range("K11").NumberFormat = "@"

$ws.Range("A11").Value = "PA"
$ws.Range("B11").Value = "PA15110"
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = "HD 300 TB"
$ws.Range("G11").Value = "44445"
$ws.Range("H11").Value = "LIDER - (9809876 01/11/25_12H) - PA"
$ws.Range("I11").Value = "01/11/25"
$ws.Range("J11").Value = "12H"
$ws.Range("K11").Value = "15/11/25"
$ws.Range("L11").Value = "DENTRO"
$ws.Range("M11").Value = ""

$ws.Range("G11").ClearFormats()
$ws.Range("I11").ClearFormats()
$ws.Range("K11").ClearFormats()
